# Trade #15 closed at 2026-02-16 21:57:44 - leadlag UP +0.000%
#
# Appends the new (still-OPEN) trade row to both the "All Trades" log and
# the strategy-specific "leadlag" sheet. The new row is cloned (copy /
# paste) from the most recent existing trade row so that text-like values
# (dates, times, etc.) keep their original text representation instead of
# being re-interpreted (e.g. "2026-02-16" auto-converting to a date
# serial) - then the fields that actually changed for this trade are
# overwritten on top.

$wb = $excel.ActiveWorkbook

$tradeNumber = 15
$tradeTime   = "21:57:44"
$entryPrice  = 68373.67999999999
$capAfter    = 100.0642567796689
$confidence  = 0.6439
$entryReason = "Coinbase leading with 0.064% move"

function Add-TradeRow {
    param(
        [string]$SheetName,
        [int]$TemplateRow
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $newRow = $TemplateRow + 1

    # Clone the previous (still OPEN) trade row so formatting / text-typed
    # cells (Date, Time, blank Exit Price, blank Exit Reason, ...) carry
    # over unchanged, then patch in the fields specific to this trade.
    $ws.Range("A" + $TemplateRow + ":O" + $TemplateRow).Copy() | Out-Null
    $ws.Range("A" + $newRow + ":O" + $newRow).PasteSpecial() | Out-Null

    $ws.Cells.Item($newRow, 1).Value  = $tradeNumber    # Trade #
    $ws.Cells.Item($newRow, 3).Value  = $tradeTime       # Time
    $ws.Cells.Item($newRow, 6).Value  = $entryPrice      # Entry Price
    $ws.Cells.Item($newRow, 11).Value = $capAfter        # Capital After
    $ws.Cells.Item($newRow, 12).Value = $confidence      # Confidence
    $ws.Cells.Item($newRow, 13).Value = $entryReason     # Entry Reason
}

Add-TradeRow "All Trades" 15
Add-TradeRow "leadlag" 14
